$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("litterChemistry")

# Clear specific cell values (keep formatting/style), matching the diff
$ws.Range("B3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("G6").ClearContents()

# Update the active selection on the sheet to D3
$ws.Activate()
$ws.Range("D3").Select()
